$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.823.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.632.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.662.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.857.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.804.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.35%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("E25").Value = "  -3.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  -3.00%  "
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.120.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.799"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").Value = "0.0₆0109"
$ws.Range("E45").Value = "  -2.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("E47").Value = "  -4.99%  "
$ws.Range("E48").Value = "  +9.26%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("E51").Value = "  -0.26%  "
